$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price (D) and volume-change (E) columns; rows 41/42 also
# swap Coin name (B) and Link (C) between MXToken and TrustWalletToken.
# A leading apostrophe forces text entry so purely-numeric-looking values
# (e.g. "213.20", "0.510") are not auto-converted to numbers by Excel,
# matching the original inline-string cell type.

$ws.Range("D2").Value = "'26.976.24"
$ws.Range("E2").Value = "'  -0.80%  "

$ws.Range("D3").Value = "'1.619.47"
$ws.Range("E3").Value = "'  -1.14%  "

$ws.Range("E4").Value = "'  +0.06%  "

$ws.Range("D5").Value = "'213.20"
$ws.Range("E5").Value = "'  -1.77%  "

$ws.Range("D6").Value = "'0.510"
$ws.Range("E6").Value = "'  -1.47%  "

$ws.Range("E7").Value = "'  +0.06%  "

$ws.Range("D8").Value = "'0.0626"
$ws.Range("E8").Value = "'  +0.10%  "

$ws.Range("E9").Value = "'  -1.60%  "

$ws.Range("D10").Value = "'19.91"
$ws.Range("E10").Value = "'  -1.09%  "

$ws.Range("D11").Value = "'0.0837"
$ws.Range("E11").Value = "'  -1.36%  "

$ws.Range("D12").Value = "'1.847.57"
$ws.Range("E12").Value = "'  -1.07%  "

$ws.Range("D13").Value = "'1.622.47"
$ws.Range("E13").Value = "'  -0.96%  "

$ws.Range("D14").Value = "'4.10"
$ws.Range("E14").Value = "'  -0.84%  "

$ws.Range("D15").Value = "'0.535"
$ws.Range("E15").Value = "'  -1.45%  "

$ws.Range("D16").Value = "'26.962.69"
$ws.Range("E16").Value = "'  -0.83%  "

$ws.Range("D17").Value = "'64.10"
$ws.Range("E17").Value = "'  -3.56%  "

$ws.Range("D18").Value = "'0.0₃0733"
$ws.Range("E18").Value = "'  -0.80%  "

$ws.Range("D19").Value = "'213.06"
$ws.Range("E19").Value = "'  -1.93%  "

$ws.Range("E20").Value = "'  +0.07%  "

$ws.Range("D21").Value = "'6.79"
$ws.Range("E21").Value = "'  -1.03%  "

$ws.Range("D22").Value = "'4.31"
$ws.Range("E22").Value = "'  -2.45%  "

$ws.Range("E23").Value = "'  -8.06%  "

$ws.Range("D24").Value = "'8.93"
$ws.Range("E24").Value = "'  -2.26%  "

$ws.Range("D25").Value = "'146.67"
$ws.Range("E25").Value = "'  -0.56%  "

$ws.Range("D26").Value = "'7.47"
$ws.Range("E26").Value = "'  +1.13%  "

$ws.Range("E27").Value = "'  +0.11%  "

$ws.Range("E28").Value = "'  -3.88%  "

$ws.Range("D29").Value = "'15.47"
$ws.Range("E29").Value = "'  -1.35%  "

$ws.Range("E30").Value = "'  -0.11%  "

$ws.Range("E31").Value = "'  -1.27%  "

$ws.Range("E32").Value = "'  -2.86%  "

$ws.Range("D33").Value = "'0.703"
$ws.Range("E33").Value = "'  +27.97%  "

$ws.Range("E34").Value = "'  -1.32%  "

$ws.Range("D35").Value = "'1.340.19"
$ws.Range("E35").Value = "'  +2.89%  "

$ws.Range("E36").Value = "'  -1.33%  "

$ws.Range("E37").Value = "'  -0.42%  "

$ws.Range("E38").Value = "'  -0.99%  "

$ws.Range("D39").Value = "'0.838"
$ws.Range("E39").Value = "'  -1.88%  "

$ws.Range("E40").Value = "'  +0.07%  "

$ws.Range("B41").Value = "'TrustWalletToken"
$ws.Range("C41").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.795"
$ws.Range("E41").Value = "'  -1.87%  "

$ws.Range("B42").Value = "'MXToken"
$ws.Range("C42").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.22"
$ws.Range("E42").Value = "'  -1.42%  "

$ws.Range("D43").Value = "'5.33"
$ws.Range("E43").Value = "'  -0.16%  "

$ws.Range("D44").Value = "'63.70"
$ws.Range("E44").Value = "'  +1.97%  "

$ws.Range("D45").Value = "'1.757.94"
$ws.Range("E45").Value = "'  -1.12%  "

$ws.Range("D46").Value = "'89.70"
$ws.Range("E46").Value = "'  -1.05%  "

$ws.Range("E47").Value = "'  +1.38%  "

$ws.Range("D48").Value = "'0.797"
$ws.Range("E48").Value = "'  +5.46%  "

$ws.Range("E49").Value = "'  +0.14%  "

$ws.Range("D50").Value = "'0.0986"
$ws.Range("E50").Value = "'  +2.90%  "

$ws.Range("D51").Value = "'7.55"
$ws.Range("E51").Value = "'  -0.87%  "
